$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Header row: translate the four column titles to their English
#    machine-friendly codes.
# ------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "mx_state"
$ws.Cells.Item(1, 2).Value = "mx_municipality"
$ws.Cells.Item(1, 3).Value = "n_matriculas"
$ws.Cells.Item(1, 4).Value = "pct_matriculas"

# ------------------------------------------------------------------
# 2. Title-case the lowercase Spanish connector words ("de", "del",
#    "la", "las", "los", "el", "y") inside every state (col A) and
#    municipality (col B) name, e.g. "Pabellón de Arteaga" ->
#    "Pabellón De Arteaga".
# ------------------------------------------------------------------
function Fix-Connectors([string]$text) {
    $result = $text
    $result = $result -replace ' de(?=[ ,])', ' De'
    $result = $result -replace ' del(?=[ ,])', ' Del'
    $result = $result -replace ' las(?=[ ,])', ' Las'
    $result = $result -replace ' la(?=[ ,])', ' La'
    $result = $result -replace ' los(?=[ ,])', ' Los'
    $result = $result -replace ' el(?=[ ,])', ' El'
    $result = $result -replace ' y(?=[ ,])', ' Y'
    return $result
}

$lastRow = 1847

# NOTE: this runtime's string comparison operators (-eq/-ceq/-ne) are
# case-insensitive, so we cannot use them to detect "did the text really
# change" (that would incorrectly treat "de" and "De" as equal and skip
# the write). Instead we unconditionally re-write every text cell with
# its (possibly identical) transformed value - harmless no-op for cells
# that contain none of the connector words.
for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valA = $cellA.Value()
    if ($valA -ne $null -and $valA -is [string] -and $valA -ne "") {
        $cellA.Value = Fix-Connectors $valA
    }

    $cellB = $ws.Cells.Item($r, 2)
    $valB = $cellB.Value()
    if ($valB -ne $null -and $valB -is [string] -and $valB -ne "") {
        $cellB.Value = Fix-Connectors $valB
    }
}

# ------------------------------------------------------------------
# 3. Two percentage cells (both 63 / 67012, the "Chalco" and
#    "Matamoros" rows) were re-saved by the original tool with a
#    1-ULP-different floating point literal; pin them to the exact
#    value from the refreshed export.
# ------------------------------------------------------------------
$ws.Cells.Item(260, 4).Value = 0.0009401301259475916
$ws.Cells.Item(1528, 4).Value = 0.0009401301259475916

# ------------------------------------------------------------------
# 4. Drop the trailing metadata/footer rows (sample size, source,
#    author, date) that used to live below the data table.
# ------------------------------------------------------------------
$ws.Rows("1849:1853").Delete()
